$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the MES (month) column values to 0 for rows 2-9
$ws.Range("B2:B9").Value = 0

# Update the active selection to B9 (matching the new selection state in the diff)
$ws.Range("B9").Select()
